$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$rows = @(
    @{ R = 367; A = 44441; B = 3; C = 27; D = 157.5722206011089 },
    @{ R = 368; A = 44442; B = 4; C = 29; D = 169.2442369419317 },
    @{ R = 369; A = 44443; B = 5; C = 33; D = 192.5882696235775 },
    @{ R = 370; A = 44444; B = 0; C = 26; D = 151.7362124306974 },
    @{ R = 371; A = 44445; B = 3; C = 25; D = 145.900204260286 },
    @{ R = 372; A = 44446; B = 0; C = 15; D = 87.54012255617158 },
    @{ R = 373; A = 44447; B = 0; C = 15; D = 87.54012255617158 },
    @{ R = 374; A = 44448; B = 1; C = 13; D = 75.8681062153487 }
)

# Copy formatting from the last existing data row (366) onto the new rows so the
# date column keeps its bordered / bold / date-formatted style (s="2") while the
# other columns keep the default (unstyled) look, exactly like all previous rows.
$ws.Range("A366:D366").Copy()
foreach ($row in $rows) {
    $ws.Range("A" + $row.R + ":D" + $row.R).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($row in $rows) {
    $ws.Cells.Item($row.R, 1).Value = $row.A
    $ws.Cells.Item($row.R, 2).Value = $row.B
    $ws.Cells.Item($row.R, 3).Value = $row.C
    $ws.Cells.Item($row.R, 4).Value = $row.D
}
